$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.847.65"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.082.66"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.30"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.42%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0789"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  +2.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.777"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.35"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "2.073.78"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "37.763.58"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.57"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.38"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.60"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.138"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.49"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.73"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.73"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.68%  "
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.41"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.20"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.74"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.19%  "
$ws.Range("D44").Value = "1.446.13"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.62%  "
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.81"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.70%  "
